# B1--and-B2-PowerPoint.pptx edit
#
# 1) The table on slide 5 gets a new table style (tableStyleId GUID change).
# 2) The deck's theme colour scheme is swapped from the "Integral / Red
#    Violet" palette to the "Office Theme / Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Change the table style on slide 5 -------------------------------
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{DA9AFAEB-A9FA-4337-AE31-15D5C09767B7}", $true)

# --- 2. Swap the theme colour scheme ------------------------------------
function ConvertTo-BGR($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$officeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$colorScheme = $p.SlideMaster.ColorScheme
for ($i = 0; $i -lt $officeColors.Length; $i++) {
    $colorScheme.Colors($i + 1).RGB = ConvertTo-BGR $officeColors[$i]
}
